$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    $s = [string]$v
    $s2 = $s -replace "_old$", "_FV2404"
    $cell.Value = $s2
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    $s = [string]$v
    $s2 = $s -replace "_new$", "_FV2410"
    $cell.Value = $s2
}

Write-Host "done renaming"

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    Write-Host "$c -> $v"
}

# Add table over the data range
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U68"), 0, 1)
$tbl.Name = "Table1"
Write-Host "table added: $($tbl.Name)"

# Freeze top row
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true
Write-Host "freeze applied"
